# Fruta / hortaliza, semanal
# A new weekly price observation is inserted as row 45 (Vega Modelo de Temuco -
# Papaya), pushing the previously existing rows 45-55 down to rows 46-56.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 45 (shifts rows 45:55 -> 46:56)
$ws.Rows("45:45").Insert()

# Populate the freshly inserted row 45 with the new weekly record
$ws.Cells.Item(45, 1).Value = 10
$ws.Cells.Item(45, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(45, 3).Value = "La Araucanía"
$ws.Cells.Item(45, 4).Value = 44466
$ws.Cells.Item(45, 5).Value = 9
$ws.Cells.Item(45, 6).Value = "Fruta"
$ws.Cells.Item(45, 7).Value = 100108
$ws.Cells.Item(45, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(45, 9).Value = 100108004
$ws.Cells.Item(45, 10).Value = "Papaya"
$ws.Cells.Item(45, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(45, 12).Value = "Primera"
$ws.Cells.Item(45, 13).Value = 180
$ws.Cells.Item(45, 14).Value = 20000
$ws.Cells.Item(45, 15).Value = 21000
$ws.Cells.Item(45, 16).Value = 20500
$ws.Cells.Item(45, 17).Value = "`$/caja 10 kilos"
$ws.Cells.Item(45, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(45, 19).Value = 2050
$ws.Cells.Item(45, 20).Value = 10
